# The deck has a single template slide (placeholder "First1"/"Last1" name
# badge text) that carries a notes page with build instructions for the
# template. Refresh the slide by duplicating it and dropping the original,
# which bumps the slide id forward and leaves the working slide with no
# notes page attached (matching the regenerated test fixture).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$newSlide = $s.Duplicate()
$s.Delete()
